$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 12.95339999999999
$ws.Range("E6").Value = 11.8111
$ws.Range("C7").Value = -11.7297
$ws.Range("E7").Value = 12.3454
$ws.Range("A8").Value = -21.1575
$ws.Range("E8").Value = 12.85980000000001
$ws.Range("E9").Value = 9.817899999999986
$ws.Range("A10").Value = -20.52489999999996
$ws.Range("E10").Value = 11.7905
$ws.Range("A12").Value = -22.46210000000004
$ws.Range("E12").Value = 12.69739999999999
$ws.Range("B13").Value = 6.477999999999998
$ws.Range("A18").Value = -22.37030000000003
$ws.Range("C20").Value = -14.93750000000001
$ws.Range("A25").Value = -22.23350000000003
